$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A2 value and apply left-aligned style
$ws.Range("A2").Value = 138158
$ws.Range("A2").HorizontalAlignment = -4131  # xlLeft

# Remove row 3 (previously A3 = 71475)
$ws.Range("A3").EntireRow.Delete()

# Update selection to match target state
$ws.Range("A3:XFD4").Select()
